$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture original values for columns D, L, M, N, O, P, Q, S, T (rows 2-19)
# before they are overwritten, since the edit is a permutation of rows.
$cols = @("D","L","M","N","O","P","Q","S","T")
$orig = @{}
for ($r = 2; $r -le 19; $r++) {
    $orig[$r] = @{}
    foreach ($c in $cols) {
        $orig[$r][$c] = $ws.Range("$c$r").Value2
    }
}

# Target row <- source row mapping derived from the diff
$mapping = @{
    2 = 16
    3 = 9
    4 = 2
    5 = 13
    6 = 8
    7 = 17
    8 = 3
    9 = 15
    10 = 4
    11 = 19
    12 = 14
    13 = 11
    14 = 18
    15 = 10
    16 = 5
    17 = 7
    18 = 6
    19 = 12
}

foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value2 = $orig[$srcRow][$c]
    }
}
